# Updated symbol list on Mon Dec 26 23:35:28 UTC 2022 with GitHub Actions
#
# Re-applies the latest crypto price/volume-label scrape onto the sheet.
# Price cells in column D are stored as literal text (they come from a
# scraped table), so a leading apostrophe is used when the new value looks
# like a number, forcing Excel to keep storing it as text (matching the
# original inlineStr/text representation) instead of silently coercing it
# to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $ws.Range($Address).Value = "'" + $Text
}

# Column D: Price
Set-TextValue "D2"  "243.39"
Set-TextValue "D4"  "5.424"
Set-TextValue "D5"  "0.05923"
Set-TextValue "D6"  "3.441"
Set-TextValue "D7"  "6.520"
Set-TextValue "D8"  "0.8089"
Set-TextValue "D9"  "0.9266"
Set-TextValue "D10" "0.1434"
Set-TextValue "D11" "0.07426"
Set-TextValue "D12" "0.03258"
Set-TextValue "D14" "0.09373"
Set-TextValue "D15" "3.869"
Set-TextValue "D16" "0.001567"
Set-TextValue "D17" "0.04678"
Set-TextValue "D18" "0.0005907"
Set-TextValue "D19" "0.005963"
Set-TextValue "D20" "0.001259"
Set-TextValue "D22" "0.00006807"
Set-TextValue "D23" "3.570"
Set-TextValue "D24" "2.134"
Set-TextValue "D26" "0.1295"
Set-TextValue "D27" "0.0002304"
Set-TextValue "D40" "0.03964"
Set-TextValue "D41" "0.006426"
Set-TextValue "D42" "0.1074"
Set-TextValue "D43" "0.003003"
Set-TextValue "D44" "0.008792"
Set-TextValue "D45" "0.00005239"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D47" "0.6708"
Set-TextValue "D48" "0.002367"
Set-TextValue "D49" "0.00002102"
Set-TextValue "D50" "0.0002002"

# Column E: Volume(1h) label -- the "Worstin24h" badge moved from row 18
# (One / ONE) to row 47 (CoinbaseStockToken / COIN).
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
